$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61748
$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -188736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1550
$ws.Range("I63").Value = 1550
$ws.Range("K63").Value = 1550
$ws.Range("M63").Value = -864
$ws.Range("H66").Value = 1550
$ws.Range("I66").Value = 1550
$ws.Range("K66").Value = 7750
$ws.Range("M66").Value = -4318
$ws.Range("H74").Value = 1409.3
$ws.Range("I74").Value = 947
$ws.Range("J74").Value = 2102.75
$ws.Range("K74").Value = 947
$ws.Range("L74").Value = 2102.75
$ws.Range("M74").Value = -73
$ws.Range("N74").Value = -3850.75
$ws.Range("H77").Value = 1409.3
$ws.Range("I77").Value = 947
$ws.Range("J77").Value = 2102.75
$ws.Range("K77").Value = 4735
$ws.Range("L77").Value = 10513.75
$ws.Range("M77").Value = -367
$ws.Range("N77").Value = -19249.75
$ws.Range("H102").Value = 9264006
$ws.Range("I102").Value = 10106007
$ws.Range("K102").Value = 10106007
$ws.Range("M102").Value = -10104385
$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2793.3547
$ws.Range("I105").Value = 2158.25
$ws.Range("J105").Value = 4970.857
$ws.Range("K105").Value = 2158.25
$ws.Range("L105").Value = 4970.857
$ws.Range("M105").Value = -411.25
$ws.Range("N105").Value = -8464.857
$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178
$ws.Range("H134").Value = 3833.3333
$ws.Range("I134").Value = 3400
$ws.Range("K134").Value = 10200
$ws.Range("M134").Value = -7665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2354
$ws.Range("J58").Value = 4614.75
$ws.Range("L58").Value = 4614.75
$ws.Range("N58").Value = -5020.75
$ws.Range("H132").Value = 1472.0714
$ws.Range("I132").Value = 1392.5
$ws.Range("J132").Value = 1949.5
$ws.Range("K132").Value = 4177.5
$ws.Range("L132").Value = 5848.5
$ws.Range("M132").Value = -1647.5
$ws.Range("N132").Value = -10908.5
$ws.Range("H136").Value = 2354
$ws.Range("J136").Value = 4614.75
$ws.Range("L136").Value = 13844.25
$ws.Range("N136").Value = -18944.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 253.23334
$ws.Range("I5").Value = 257.88
$ws.Range("J5").Value = 230
$ws.Range("K5").Value = 773.64
$ws.Range("L5").Value = 690
$ws.Range("M5").Value = -661.64
$ws.Range("N5").Value = -914
$ws.Range("H68").Value = 390
$ws.Range("I68").Value = 390
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1170
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -359
$ws.Range("N68").ClearContents()
$ws.Range("H70").Value = 799.5
$ws.Range("I70").Value = 799.5
$ws.Range("K70").Value = 2398.5
$ws.Range("M70").Value = -2083.5
$ws.Range("H71").Value = 390
$ws.Range("I71").Value = 390
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 3510
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 546
$ws.Range("N71").ClearContents()
$ws.Range("H73").Value = 799.5
$ws.Range("I73").Value = 799.5
$ws.Range("K73").Value = 2398.5
$ws.Range("M73").Value = -1306.5
$ws.Range("H75").Value = 1234
$ws.Range("J75").Value = 1234
$ws.Range("L75").Value = 3702
$ws.Range("N75").Value = -5698
$ws.Range("H78").Value = 1234
$ws.Range("J78").Value = 1234
$ws.Range("L78").Value = 11106
$ws.Range("N78").Value = -21090
$ws.Range("H82").Value = 11013
$ws.Range("I82").Value = 11013
$ws.Range("K82").Value = 33039
$ws.Range("M82").Value = -32633
$ws.Range("H85").Value = 11013
$ws.Range("I85").Value = 11013
$ws.Range("K85").Value = 33039
$ws.Range("M85").Value = -31635
$ws.Range("H87").Value = 1000
$ws.Range("I87").Value = 1000
$ws.Range("K87").Value = 3000
$ws.Range("M87").Value = -1752
$ws.Range("H90").Value = 1000
$ws.Range("I90").Value = 1000
$ws.Range("K90").Value = 9000
$ws.Range("M90").Value = -2760
$ws.Range("H103").Value = 204.8
$ws.Range("I103").Value = 156
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 468
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = 411
$ws.Range("N103").Value = -2958
$ws.Range("H107").Value = 167111.5
$ws.Range("J107").Value = 250392.25
$ws.Range("L107").Value = 751176.75
$ws.Range("N107").Value = -755016.75
$ws.Range("H109").Value = 1086.875
$ws.Range("J109").Value = 1899.5
$ws.Range("L109").Value = 5698.5
$ws.Range("N109").Value = -7778.5
$ws.Range("H114").Value = 1061
$ws.Range("I114").Value = 842.6667
$ws.Range("J114").Value = 1224.75
$ws.Range("K114").Value = 2528.0001
$ws.Range("L114").Value = 3674.25
$ws.Range("M114").Value = 725.9998999999998
$ws.Range("N114").Value = -10182.25
$ws.Range("H122").Value = 230.71428
$ws.Range("I122").Value = 219.16667
$ws.Range("K122").Value = 1972.50003
$ws.Range("M122").Value = 477.4999699999998
$ws.Range("H132").Value = 2823.1667
$ws.Range("J132").Value = 2600
$ws.Range("L132").Value = 23400
$ws.Range("N132").Value = -28460
$ws.Range("H135").Value = 253.23334
$ws.Range("I135").Value = 257.88
$ws.Range("J135").Value = 230
$ws.Range("K135").Value = 2320.92
$ws.Range("L135").Value = 2070
$ws.Range("M135").Value = 214.0799999999999
$ws.Range("N135").Value = -7140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 30000.5
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 30000.5
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 30000.5
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -31262.5
$ws.Range("H70").Value = 20198
$ws.Range("I70").Value = 12998.5
$ws.Range("J70").Value = 24997.666
$ws.Range("K70").Value = 12998.5
$ws.Range("L70").Value = 24997.666
$ws.Range("M70").Value = -12728.5
$ws.Range("N70").Value = -25537.666
$ws.Range("H73").Value = 20198
$ws.Range("I73").Value = 12998.5
$ws.Range("J73").Value = 24997.666
$ws.Range("K73").Value = 12998.5
$ws.Range("L73").Value = 24997.666
$ws.Range("M73").Value = -12062.5
$ws.Range("N73").Value = -26869.666
$ws.Range("H113").Value = 13916.667
$ws.Range("I113").Value = 1250
$ws.Range("K113").Value = 1250
$ws.Range("M113").Value = 920
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H16").Value = 1472.5
$ws.Range("I16").Value = 1300
$ws.Range("K16").Value = 1300
$ws.Range("M16").Value = -1130
$ws.Range("H40").Value = 3178.7778
$ws.Range("I40").Value = 3135.7917
$ws.Range("K40").Value = 3135.7917
$ws.Range("M40").Value = -2999.7917
$ws.Range("H116").Value = 20000
$ws.Range("J116").Value = 20000
$ws.Range("L116").Value = 20000
$ws.Range("N116").Value = -29178
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 10000
$ws.Range("J86").Value = 10000
$ws.Range("L86").Value = 10000
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 10000
$ws.Range("J89").Value = 10000
$ws.Range("L89").Value = 50000
$ws.Range("N89").Value = -61232
$ws.Range("H107").Value = 25496
$ws.Range("I107").Value = 25496
$ws.Range("K107").Value = 76488
$ws.Range("M107").Value = -74568
$ws.Range("H110").Value = 83000
$ws.Range("I110").Value = 83000
$ws.Range("K110").Value = 83000
$ws.Range("M110").Value = -78910
$ws.Range("H116").Value = 44000
$ws.Range("J116").Value = 44000
$ws.Range("L116").Value = 44000
$ws.Range("N116").Value = -53178
$ws.Range("H117").Value = 75409
$ws.Range("J117").Value = 75409
$ws.Range("L117").Value = 75409
$ws.Range("N117").Value = -84587
$ws.Range("H122").Value = 2598.3333
$ws.Range("I122").Value = 2398.5625
$ws.Range("K122").Value = 7195.6875
$ws.Range("M122").Value = -4745.6875
